# Lab Exam 03 grading workbook - "kalyankar to pusapati done"
# Fill in the "Points for grading" (column E) for the CustomerMapping Class
# Generic section (rows 3-6) and the Customer Class section (rows 10-14),
# matching the "Total Points" already entered in column D for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generic section (S.No 1-4)
$ws.Range("E3").Value  = 1
$ws.Range("E4").Value  = 2
$ws.Range("E5").Value  = 2
$ws.Range("E6").Value  = 2

# Customer Class section (S.No 5-9)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the current selection to E15 (the running total for the section
# just graded), matching where the grader's cursor ended up.
$ws.Range("E15").Select()
